$wb = $excel.ActiveWorkbook

# --- ARM sheet: clear stale price/profit data (columns H:N) for rows 121-141,
#     except row 136 which keeps its values ---
$wsArm = $wb.Worksheets.Item("ARM")
$wsArm.Range("H121:N135").ClearContents()
$wsArm.Range("H137:N141").ClearContents()

# --- CUL sheet: clear stale price/profit data (columns H:N) for rows 120-141,
#     except row 135 which keeps its values ---
$wsCul = $wb.Worksheets.Item("CUL")
$wsCul.Range("H120:N134").ClearContents()
$wsCul.Range("H136:N141").ClearContents()

# --- LTW sheet: refresh price/profit data for rows 61 and 113 ---
$wsLtw = $wb.Worksheets.Item("LTW")

$wsLtw.Range("H61").Value = 1189.6666
$wsLtw.Range("I61").Value = 941.4
$wsLtw.Range("J61").Value = 1500
$wsLtw.Range("K61").Value = 941.4
$wsLtw.Range("L61").Value = 1500
$wsLtw.Range("M61").Value = -739.4
$wsLtw.Range("N61").Value = -1904

$wsLtw.Range("H113").Value = 1189.6666
$wsLtw.Range("I113").Value = 941.4
$wsLtw.Range("J113").Value = 1500
$wsLtw.Range("K113").Value = 941.4
$wsLtw.Range("L113").Value = 1500
$wsLtw.Range("M113").Value = 1228.6
$wsLtw.Range("N113").Value = -5840
